$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(5,1).NumberFormat = "@"
$ws.Cells.Item(5,1).Value = "2025-09-21"

$ws.Cells.Item(5,2).Value = "Pick 3"

$ws.Cells.Item(5,3).NumberFormat = "@"
$ws.Cells.Item(5,3).Value = "250921"

$ws.Cells.Item(5,4).Value = "3-2-8"

$ws.Cells.Item(5,5).Value = "2025-09-21T21:35:39.672+04:00"
